$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "On this page" summary cell (A7) with the new date/time.
$ws.Range("A7").Value = "On this pageCurrent school and early childhood service, TAFE closures and relocations:Bus service cancellations or alterationsCurrent school and early childhood service, TAFE closures and relocations for Thursday 3 September, (as at 10:30am, 3 September)South-Eastern Victoria RegionEarly childhood services"

# 2. Insert a new early-childhood-service closure entry before the old row 177
#    ("Camp Australia - Haileybury City Campus OSHC WEST MELBOURNE"), pushing
#    it and everything below down by one row.
$ws.Rows(177).Insert()
$ws.Range("A177").Value = "li: Camp Australia - Flemington Primary School OSHC FLEMINGTON"

# 3. Insert another new entry before the old row 190 ("Hopetoun Child Care
#    Service HOPETOUN"), which after step 2 now sits at row 191.
$ws.Rows(191).Insert()
$ws.Range("A191").Value = "li: Gowrie Victoria Clare Court YARRAVILLE"

# 4. Update the "Last Update" footer line (now at row 223 after the two
#    insertions above) to reflect the new publish date.
$ws.Range("A223").Value = "li: Last Update: 03 September 2020"
